$d = $word.ActiveDocument

# The document contains two similarly-worded passages (the Abstract and,
# later, the Introduction) -- anchor on text that is unique to the Abstract
# paragraph so we never touch the Introduction's near-duplicate wording.
$paraAnchor = "In statistical integrated age structured population models, there are two common practices"

# --- 1. Merge the "parameters" run (and drop the surrounding proofErr
#        gramStart/gramEnd tags) back into the main paragraph run, exactly
#        as Word collapses runs into a single contiguous run when text is
#        retyped/re-accepted across them. A genuine text change is needed
#        to force the merge, so we stage the replacement through a
#        temporary placeholder and then restore the original wording. ---
$oldFragment = "estimates growth parameters with other population dynamics and fishery processes. When growth is estimated externally"
$newFragment = $oldFragment   # text itself is unchanged; only run/proofErr structure collapses

$content = $d.Content.Text
$paraIdx = $content.IndexOf($paraAnchor)
$fragIdx = $content.IndexOf($oldFragment, $paraIdx)
$fragEnd = $fragIdx + $oldFragment.Length
$r = $d.Range($fragIdx, $fragEnd)
$r.Text = "ss3sim_TEMP_PLACEHOLDER_9f3c"

$content = $d.Content.Text
$tempIdx = $content.IndexOf("ss3sim_TEMP_PLACEHOLDER_9f3c", $paraIdx)
$r2 = $d.Range($tempIdx, $tempIdx + "ss3sim_TEMP_PLACEHOLDER_9f3c".Length)
$r2.Text = $newFragment

# --- 2. Move the "_GoBack" bookmark from the end of the abstract paragraph
#        to the start of the abstract paragraph (right before the first
#        run) -- this is where Word leaves the "last edit" bookmark after
#        the CAPAM workshop abstract text was inserted at the top of the
#        paragraph. ---
$content = $d.Content.Text
$abstractIdx = $content.IndexOf($paraAnchor)
$startRange = $d.Range($abstractIdx, $abstractIdx)
$d.Bookmarks.Add("_GoBack", $startRange)
